# Applies the "LinuxForHealth" re-branding edit to the
# StructureDefinition-cost-sharing-reduction-variant workbook:
#   - Metadata sheet: URL, Version, Date and Publisher values are refreshed.
#   - Elements sheet: the combined ele-1/ext-1 invariant text is moved off the
#     top-level "Extension" row (row 2) and onto the "Extension.extension"
#     row (row 4) only, in the "Constraint(s)" column (AI).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cost-sharing-reduction-variant"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Column AI = "Constraint(s)".
# Row 2 = "Extension" (base row)      -> drop the constraint note.
# Row 4 = "Extension.extension" row   -> gains the constraint note.
$constraintText = $elements.Range("AI2").Value2

$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraintText
